# Applies the "use IAPS pictures, make them big, and with constant size" edit:
#  - cue_values!B2:B37 (IAPS stimulus filenames) collapse to a small constant
#    set of IAPS picture filenames, mostly the single "1030.JPG" image.
#  - cue_validities!B1:B6 text colour becomes explicit black (was theme colour).
#  - Assorted row-height bumps (18.75 -> 19.5) on cue_validities / strat_pred.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("cue_values")
$ws2 = $wb.Worksheets.Item("cue_validities")
$ws3 = $wb.Worksheets.Item("strat_pred")

# --- cue_values: column B (IAPSslide) trial stimulus names -----------------
# Row 2 / 3 get their own distinct filenames, every remaining row (4-37)
# reuses the same constant-size picture.
$ws1.Cells.Item(2, 2).Value2 = "1710.JPG"
$ws1.Cells.Item(3, 2).Value2 = "1022.JPG"
for ($r = 4; $r -le 37; $r++) {
    $ws1.Cells.Item($r, 2).Value2 = "1030.JPG"
}

# --- cue_validities: make the label font explicit black --------------------
$ws2.Range("B1:B6").Font.Color = 0

# --- row-height bumps (18.75pt -> 19.5pt) -----------------------------------
for ($r = 1; $r -le 6; $r++) {
    $ws2.Rows.Item($r).RowHeight = 19.5
}
for ($r = 1; $r -le 37; $r++) {
    $ws3.Rows.Item($r).RowHeight = 19.5
}
